$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking text (prices with "." as thousands
# separators, e.g. "64.609.86") that must stay literal text. Force the
# whole data range to Text format before writing so Excel does not
# auto-convert values such as "180.42" or "11.24" into real numbers,
# then restore the Normal style so no stray number-format/style is left
# behind on the cells.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "64.609.86"
$ws.Range("E2").Value = "  -3.60%  "
$ws.Range("D3").Value = "3.313.98"
$ws.Range("E3").Value = "  -4.60%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "180.42"
$ws.Range("E5").Value = "  -8.97%  "
$ws.Range("D6").Value = "529.15"
$ws.Range("E6").Value = "  -3.21%  "
$ws.Range("D7").Value = "0.605"
$ws.Range("E7").Value = "  +0.43%  "
$ws.Range("D8").Value = "3.307.09"
$ws.Range("E8").Value = "  -4.47%  "
$ws.Range("E9").Value = "  +0.04%  "
$ws.Range("D10").Value = "0.612"
$ws.Range("E10").Value = "  -5.29%  "
$ws.Range("D11").Value = "59.28"
$ws.Range("E11").Value = "  -3.82%  "
$ws.Range("E12").Value = "  -5.46%  "
$ws.Range("E13").Value = "  -2.36%  "
$ws.Range("E14").Value = "  -5.57%  "
$ws.Range("D15").Value = "3.844.56"
$ws.Range("E15").Value = "  -4.81%  "
$ws.Range("D16").Value = "3.315.43"
$ws.Range("E16").Value = "  -4.99%  "
$ws.Range("D17").Value = "0.117"
$ws.Range("E17").Value = "  -4.73%  "
$ws.Range("D18").Value = "64.583.32"
$ws.Range("E18").Value = "  -3.26%  "
$ws.Range("D19").Value = "17.61"
$ws.Range("E19").Value = "  -2.67%  "
$ws.Range("E20").Value = "  -3.94%  "
$ws.Range("D21").Value = "0.963"
$ws.Range("E21").Value = "  -4.74%  "
$ws.Range("D22").Value = "376.11"
$ws.Range("E22").Value = "  -2.43%  "
$ws.Range("D23").Value = "3.82"
$ws.Range("E23").Value = "  -3.29%  "
$ws.Range("B24").Value = "RenderToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D24").Value = "11.24"
$ws.Range("E24").Value = "  -4.65%  "
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").Value = "80.95"
$ws.Range("E25").Value = "  -1.02%  "
$ws.Range("E26").Value = "  +3.02%  "
$ws.Range("E27").Value = "  -0.71%  "
$ws.Range("E28").Value = "  -2.50%  "
$ws.Range("D29").Value = "11.59"
$ws.Range("E29").Value = "  -3.69%  "
$ws.Range("D30").Value = "8.42"
$ws.Range("E30").Value = "  -3.20%  "
$ws.Range("D31").Value = "29.10"
$ws.Range("E31").Value = "  -5.26%  "
$ws.Range("D32").Value = "656.33"
$ws.Range("E32").Value = "  -2.18%  "
$ws.Range("D33").Value = "6.70"
$ws.Range("E33").Value = "  -2.42%  "
$ws.Range("D34").Value = "11.35"
$ws.Range("E34").Value = "  -1.92%  "
$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").Value = "59.66"
$ws.Range("E35").Value = "  -5.90%  "
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").Value = "0.106"
$ws.Range("E36").Value = "  -2.49%  "
$ws.Range("D38").Value = "0.396"
$ws.Range("E38").Value = "  +0.96%  "
$ws.Range("D39").Value = "36.88"
$ws.Range("E39").Value = "  -3.12%  "
$ws.Range("E40").Value = "  +0.10%  "
$ws.Range("D41").Value = "0.0₃0705"
$ws.Range("E41").Value = "  +6.25%  "
$ws.Range("E42").Value = "  -1.23%  "
$ws.Range("D43").Value = "2.913.29"
$ws.Range("E43").Value = "  -4.20%  "
$ws.Range("D44").Value = "2.54"
$ws.Range("E44").Value = "  +3.01%  "
$ws.Range("E45").Value = "  -7.68%  "
$ws.Range("D46").Value = "0.0401"
$ws.Range("E46").Value = "  +2.62%  "
$ws.Range("E47").Value = "  -4.10%  "
$ws.Range("D48").Value = "3.12"
$ws.Range("E48").Value = "  +9.68%  "
$ws.Range("D49").Value = "2.84"
$ws.Range("E49").Value = "  +10.02%  "
$ws.Range("E50").Value = "  +1.56%  "
$ws.Range("D51").Value = "2.56"
$ws.Range("E51").Value = "  -5.02%  "

$ws.Range("D2:D51").Style = "Normal"
